$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by the type
# inference get a temporary "@" (Text) number format so the stored
# value stays a string, then the style is reset to Normal so no
# extra formatting is left behind on the cell.
$textForce = @(
    @{Cell='D2'; Value='70.130.22'}
    @{Cell='D3'; Value='3.585.76'}
    @{Cell='D5'; Value='578.84'}
    @{Cell='D6'; Value='188.75'}
    @{Cell='D8'; Value='3.581.06'}
    @{Cell='D11'; Value='0.659'}
    @{Cell='D12'; Value='56.08'}
    @{Cell='D14'; Value='9.60'}
    @{Cell='D15'; Value='4.157.33'}
    @{Cell='D16'; Value='19.84'}
    @{Cell='D17'; Value='3.582.14'}
    @{Cell='D18'; Value='70.036.91'}
    @{Cell='D19'; Value='12.59'}
    @{Cell='D22'; Value='474.15'}
    @{Cell='D23'; Value='18.94'}
    @{Cell='D24'; Value='5.07'}
    @{Cell='D26'; Value='88.78'}
    @{Cell='D28'; Value='11.07'}
    @{Cell='D29'; Value='9.33'}
    @{Cell='D30'; Value='32.14'}
    @{Cell='D31'; Value='7.71'}
    @{Cell='D33'; Value='12.13'}
    @{Cell='D34'; Value='65.97'}
    @{Cell='D35'; Value='586.19'}
    @{Cell='D36'; Value='39.03'}
    @{Cell='D38'; Value='0.0₃0796'}
    @{Cell='D39'; Value='0.395'}
    @{Cell='D40'; Value='0.139'}
    @{Cell='D41'; Value='3.51'}
    @{Cell='B42'; Value='dogwifhat'}
    @{Cell='C42'; Value='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'}
    @{Cell='D42'; Value='3.22'}
    @{Cell='B43'; Value='Fetch.AI'}
    @{Cell='C43'; Value='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'}
    @{Cell='D43'; Value='2.89'}
    @{Cell='D44'; Value='3.235.46'}
    @{Cell='D47'; Value='9.59'}
    @{Cell='D48'; Value='3.33'}
)

foreach ($item in $textForce) {
    $r = $ws.Range($item.Cell)
    $r.NumberFormat = "@"
    $r.Value = $item.Value
    $r.Style = "Normal"
}

$plainUpdates = @{
    'E3' = '  -1.02%  '
    'E4' = '  +0.07%  '
    'E6' = '  -2.30%  '
    'E7' = '  -2.50%  '
    'E8' = '  -0.18%  '
    'E9' = '  +0.04%  '
    'E10' = '  -1.93%  '
    'E11' = '  -0.60%  '
    'E12' = '  -3.23%  '
    'E13' = '  +1.59%  '
    'E14' = '  -1.74%  '
    'E15' = '  -1.19%  '
    'E16' = '  +2.49%  '
    'E17' = '  -1.28%  '
    'E18' = '  -0.34%  '
    'E19' = '  -0.45%  '
    'E20' = '  +0.34%  '
    'E21' = '  -1.15%  '
    'E22' = '  -4.48%  '
    'E23' = '  +12.92%  '
    'E24' = '  -8.66%  '
    'E25' = '  -2.18%  '
    'E26' = '  -2.27%  '
    'E27' = '  -2.58%  '
    'E28' = '  -1.42%  '
    'E29' = '  -0.82%  '
    'E30' = '  -0.64%  '
    'E31' = '  +2.02%  '
    'E32' = '  +3.26%  '
    'E33' = '  -0.94%  '
    'E34' = '  +0.99%  '
    'E35' = '  -5.14%  '
    'E36' = '  +2.71%  '
    'E37' = '  -0.03%  '
    'E38' = '  -4.67%  '
    'E39' = '  -1.62%  '
    'E40' = '  -6.59%  '
    'E41' = '  -5.70%  '
    'E42' = '  +15.92%  '
    'E43' = '  +7.57%  '
    'E44' = '  -3.09%  '
    'E46' = '  -1.57%  '
    'E47' = '  +5.59%  '
    'E48' = '  +0.64%  '
    'E49' = '  -0.67%  '
    'E50' = '  -0.17%  '
    'E51' = '  -4.28%  '
}

foreach ($cell in $plainUpdates.Keys) {
    $ws.Range($cell).Value = $plainUpdates[$cell]
}
